# Apply the daily-update edits to the multiplication practice sheet:
# - Update the date header
# - Update each "AAA x B=" arithmetic problem to its new values

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2025-11-05 Wednesday"; New = "2025-11-06 Thursday" },
    @{ Old = "683×2=";  New = "887×9=" },
    @{ Old = "669×9=";  New = "336×9=" },
    @{ Old = "324×4=";  New = "741×4=" },
    @{ Old = "911×3=";  New = "554×8=" },
    @{ Old = "423×4=";  New = "241×7=" },
    @{ Old = "564×8=";  New = "852×2=" },
    @{ Old = "269×3=";  New = "162×2=" },
    @{ Old = "389×3=";  New = "819×5=" },
    @{ Old = "968×4=";  New = "205×5=" },
    @{ Old = "825×6=";  New = "746×6=" },
    @{ Old = "302×6=";  New = "858×7=" },
    @{ Old = "441×2=";  New = "367×4=" },
    @{ Old = "501×7=";  New = "762×5=" },
    @{ Old = "625×2=";  New = "677×7=" },
    @{ Old = "475×5=";  New = "473×7=" },
    @{ Old = "844×9=";  New = "929×5=" },
    @{ Old = "878×7=";  New = "356×8=" },
    @{ Old = "321×8=";  New = "315×9=" },
    @{ Old = "137×5=";  New = "283×2=" },
    @{ Old = "415×7=";  New = "881×6=" },
    @{ Old = "742×4=";  New = "929×7=" },
    @{ Old = "675×4=";  New = "701×9=" },
    @{ Old = "498×9=";  New = "754×3=" },
    @{ Old = "115×5=";  New = "673×8=" },
    @{ Old = "656×5=";  New = "220×3=" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $r.New, 2)
}

$d.Save()
